$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block in rows 166-206 is a weekly price history for
# "Femacal de La Calera - Ciboulette" that only varies by column
# D (Fecha / date), J (Volumen) and, on one special low-volume
# observation, K/L/M (Precio minimo/maximo/promedio) and P (Precio
# $/Kg). A new weekly observation was appended, which shifts every
# existing observation down by one row: row 166 receives the new
# date (44511) while keeping its other figures, rows 167-206 each
# take on the D/J/K/L/M/P values that used to belong to the row
# right above them, and the values that used to be in the last row
# (206) move into a brand-new row 207.

# Step 1: materialise new row 207 as a full copy of row 206 (its
# current / original content), setting each cell explicitly so the
# date cell reuses the existing date style/number format instead of
# Excel inventing a new one via a generic paste.
$ws.Range("A207").Value = $ws.Range("A206").Value2
$ws.Range("B207").Value = $ws.Range("B206").Value2
$ws.Range("C207").Value = $ws.Range("C206").Value2
$ws.Range("D207").Value = $ws.Range("D206").Value2
$ws.Range("D207").NumberFormat = $ws.Range("D206").NumberFormat
$ws.Range("E207").Value = $ws.Range("E206").Value2
$ws.Range("F207").Value = $ws.Range("F206").Value2
$ws.Range("G207").Value = $ws.Range("G206").Value2
$ws.Range("H207").Value = $ws.Range("H206").Value2
$ws.Range("I207").Value = $ws.Range("I206").Value2
$ws.Range("J207").Value = $ws.Range("J206").Value2
$ws.Range("K207").Value = $ws.Range("K206").Value2
$ws.Range("L207").Value = $ws.Range("L206").Value2
$ws.Range("M207").Value = $ws.Range("M206").Value2
$ws.Range("N207").Value = $ws.Range("N206").Value2
$ws.Range("O207").Value = $ws.Range("O206").Value2
$ws.Range("P207").Value = $ws.Range("P206").Value2
$ws.Range("Q207").Value = $ws.Range("Q206").Value2
$ws.Range("R207").Value = $ws.Range("R206").Value2

# Step 2: shift D/J/K/L/M/P down by one row for rows 167..206,
# i.e. row N takes the (still-original) values from row N-1.
# Walking from the bottom (206) up to the top (167) guarantees each
# source row is read before it is itself overwritten.
for ($n = 206; $n -ge 167; $n--) {
    $src = $n - 1
    $ws.Range("D$n").Value = $ws.Range("D$src").Value2
    $ws.Range("J$n").Value = $ws.Range("J$src").Value2
    $ws.Range("K$n").Value = $ws.Range("K$src").Value2
    $ws.Range("L$n").Value = $ws.Range("L$src").Value2
    $ws.Range("M$n").Value = $ws.Range("M$src").Value2
    $ws.Range("P$n").Value = $ws.Range("P$src").Value2
}

# Step 3: row 166 becomes the newest observation, with a brand new
# date; its Volumen/Precio figures stay the same as before.
$ws.Range("D166").Value = 44511
